$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Update the name/lastname in row 4 (Esteban Quito -> Harry Plotter)
$ws.Range("A4").Value = "Harry"
$ws.Range("B4").Value = "Plotter"

# Update the active selection to match the final state
$ws.Range("B5").Select()
